$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell holds plain text in the source workbook (t="inlineStr").
# Excel auto-converts number-looking strings ("307.38", "1.011", ...) to
# real numbers on assignment, so we briefly force Text format, assign the
# literal string, then clear the format delta back off so no stray style
# index is left on the cell (matches the un-styled cells in the source).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '26.817.29'
Set-TextValue 'E2' '  -2.48%  '
Set-TextValue 'D3' '1.779.85'
Set-TextValue 'E3' '  -2.77%  '
Set-TextValue 'D4' '1.011'
Set-TextValue 'E4' '  +0.85%  '
Set-TextValue 'E5' '  +0.62%  '
Set-TextValue 'D6' '307.38'
Set-TextValue 'E6' '  -1.89%  '
Set-TextValue 'D7' '0.4212'
Set-TextValue 'E7' '  -2.08%  '
Set-TextValue 'D8' '0.3598'
Set-TextValue 'E8' '  -1.85%  '
Set-TextValue 'D9' '0.07153'
Set-TextValue 'E9' '  -1.57%  '
Set-TextValue 'D10' '0.8340'
Set-TextValue 'E10' '  -3.80%  '
Set-TextValue 'D11' '20.19'
Set-TextValue 'D12' '1.913.15'
Set-TextValue 'E12' '  +2.37%  '
Set-TextValue 'D13' '5.232'
Set-TextValue 'E13' '  -3.17%  '
Set-TextValue 'D14' '6.318'
Set-TextValue 'E14' '  -3.18%  '
Set-TextValue 'D15' '0.06805'
Set-TextValue 'E15' '  -1.77%  '
Set-TextValue 'E16' '  +0.83%  '
Set-TextValue 'D17' '78.99'
Set-TextValue 'E17' '  -2.01%  '
Set-TextValue 'D18' '0.000008662'
Set-TextValue 'E18' '  -2.36%  '
Set-TextValue 'D19' '1.007'
Set-TextValue 'E19' '  +0.60%  '
Set-TextValue 'D20' '14.84'
Set-TextValue 'E20' '  -3.54%  '
Set-TextValue 'D21' '26.948.95'
Set-TextValue 'E21' '  -2.37%  '
Set-TextValue 'D22' '5.014'
Set-TextValue 'E22' '  -2.35%  '
Set-TextValue 'D23' '11.00'
Set-TextValue 'E23' '  +1.75%  '
Set-TextValue 'D24' '2.034.03'
Set-TextValue 'E24' '  -4.07%  '
Set-TextValue 'D25' '1.915'
Set-TextValue 'E25' '  -3.23%  '
Set-TextValue 'D26' '152.83'
Set-TextValue 'E26' '  -1.03%  '
Set-TextValue 'D27' '18.07'
Set-TextValue 'E27' '  -3.85%  '
Set-TextValue 'B28' 'BitcoinCash'
Set-TextValue 'C28' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D28' '114.26'
Set-TextValue 'E28' '  +0.24%  '
Set-TextValue 'B29' 'InternetComputer(DFINITY)'
Set-TextValue 'C29' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D29' '5.009'
Set-TextValue 'E29' '  -2.52%  '
Set-TextValue 'D30' '1.620'
Set-TextValue 'E30' '  -11.25%  '
Set-TextValue 'D31' '0.08927'
Set-TextValue 'E31' '  +0.50%  '
Set-TextValue 'D32' '2.843'
Set-TextValue 'E32' '  -4.79%  '
Set-TextValue 'D33' '0.7094'
Set-TextValue 'E33' '  -5.65%  '
Set-TextValue 'D34' '4.291'
Set-TextValue 'E34' '  -5.58%  '
Set-TextValue 'D35' '1.085'
Set-TextValue 'E35' '  -4.26%  '
Set-TextValue 'E36' '  +0.64%  '
Set-TextValue 'E37' '  -1.15%  '
Set-TextValue 'D38' '0.01883'
Set-TextValue 'E38' '  -2.69%  '
Set-TextValue 'D39' '0.05068'
Set-TextValue 'E39' '  -4.71%  '
Set-TextValue 'D40' '0.4898'
Set-TextValue 'E40' '  -3.43%  '
Set-TextValue 'D41' '0.1603'
Set-TextValue 'E41' '  -3.74%  '
Set-TextValue 'D42' '2.510'
Set-TextValue 'E42' '  -10.38%  '
Set-TextValue 'D43' '5.961'
Set-TextValue 'E43' '  -9.62%  '
Set-TextValue 'D44' '7.879'
Set-TextValue 'E44' '  -5.94%  '
Set-TextValue 'D46' '104.02'
Set-TextValue 'E46' '  -1.83%  '
Set-TextValue 'D47' '10.09'
Set-TextValue 'E47' '  -3.81%  '
Set-TextValue 'D48' '0.06243'
Set-TextValue 'E48' '  -3.94%  '
Set-TextValue 'D49' '0.4445'
Set-TextValue 'E49' '  -5.04%  '
Set-TextValue 'D50' '1.567'
Set-TextValue 'E50' '  -2.72%  '
Set-TextValue 'D51' '1.694'
Set-TextValue 'E51' '  -2.25%  '
